$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (RF)
$ws.Range("B3").Value = 0.243
$ws.Range("C3").Value = -0.092
$ws.Range("D3").Value = 0.48
$ws.Range("E3").Value = 0.693
$ws.Range("F3").Value = 0.769
$ws.Range("G3").Value = 0.545

# Row 4 (NN)
$ws.Range("B4").Value = 0.258
$ws.Range("C4").Value = -0.07
$ws.Range("D4").Value = 0.471
$ws.Range("E4").Value = 0.686
$ws.Range("F4").Value = 0.751
$ws.Range("G4").Value = 0.545

# Row 5 (RNN)
$ws.Range("B5").Value = 0.024
$ws.Range("C5").Value = -0.206
$ws.Range("D5").Value = 0.556
$ws.Range("E5").Value = 0.746
$ws.Range("F5").Value = 0.686
$ws.Range("G5").Value = 0.336

# Row 6 (Ensemble)
$ws.Range("B6").Value = 0.296
$ws.Range("C6").Value = -0.015
$ws.Range("D6").Value = 0.447
$ws.Range("E6").Value = 0.669
$ws.Range("F6").Value = 0.679
$ws.Range("G6").Value = 0.546
